$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.234.63"
$ws.Range("E2").Value = "  -0.49%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.830.25"
$ws.Range("E3").Value = "  -0.67%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.37"
$ws.Range("E5").Value = "  -1.21%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6098"
$ws.Range("E6").Value = "  -3.33%  "

$ws.Range("E7").Value = "  +0.21%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07117"
$ws.Range("E8").Value = "  -4.79%  "

$ws.Range("E9").Value = "  -2.64%  "

$ws.Range("E10").Value = "  -4.15%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07646"
$ws.Range("E11").Value = "  -1.20%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.842.40"
$ws.Range("E12").Value = "  +0.02%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.819"
$ws.Range("E13").Value = "  -3.43%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6384"
$ws.Range("E14").Value = "  -6.02%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000009969"
$ws.Range("E15").Value = "  -2.55%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.066.54"
$ws.Range("E16").Value = "  -1.23%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "79.79"
$ws.Range("E17").Value = "  -2.86%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.988"
$ws.Range("E18").Value = "  -4.81%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "29.220.26"
$ws.Range("E19").Value = "  -0.49%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "230.37"
$ws.Range("E20").Value = "  +0.43%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.83"
$ws.Range("E21").Value = "  -4.15%  "

$ws.Range("E22").Value = "  +0.17%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.050"
$ws.Range("E23").Value = "  -4.92%  "

$ws.Range("E24").Value = "  +0.19%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.50"
$ws.Range("E25").Value = "  -1.98%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.101"
$ws.Range("E26").Value = "  -4.73%  "

$ws.Range("E27").Value = "  -4.09%  "

$ws.Range("E28").Value = "  -3.82%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06803"
$ws.Range("E29").Value = "  +3.44%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.481"
$ws.Range("E30").Value = "  +2.79%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.457"
$ws.Range("E31").Value = "  -2.20%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.849"
$ws.Range("E32").Value = "  -5.15%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.830"
$ws.Range("E33").Value = "  -5.95%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.128"
$ws.Range("E34").Value = "  -1.16%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.737"
$ws.Range("E35").Value = "  -5.70%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6593"
$ws.Range("E36").Value = "  -5.80%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.554"
$ws.Range("E37").Value = "  -0.92%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.233.69"
$ws.Range("E38").Value = "  -1.25%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01766"
$ws.Range("E40").Value = "  -4.85%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.606"
$ws.Range("E41").Value = "  -2.98%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9343"
$ws.Range("E42").Value = "  +0.04%  "

$ws.Range("E43").Value = "  +0.17%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.985.74"
$ws.Range("E44").Value = "  -0.47%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "100.93"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "63.59"
$ws.Range("E46").Value = "  -3.01%  "

$ws.Range("E47").Value = "  -1.96%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.632"
$ws.Range("E48").Value = "  -5.01%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.568"
$ws.Range("E49").Value = "  -5.49%  "

# Row 50/51: Algorand/Aptos swap with updated values
$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.553"
$ws.Range("E50").Value = "  -7.30%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1089"
$ws.Range("E51").Value = "  -5.05%  "
